$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 27) below the existing last row (26)
$ws.Range("A27").Value = "E_PRZNTPKT"
$ws.Range("B27").Value = "Prozentpunkte"
$ws.Range("C27").Value = "Percentage points"

# Copy the formatting (style) from the last existing data row onto the new row
$ws.Range("A26:C26").Copy()
$ws.Range("A27:C27").PasteSpecial(-4122)
